$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Update existing data rows (2-8) with re-shuffled / changed values ----

# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "IB00VJ993_B"
$ws.Range("C2").Value = "Top       "
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = "SMT"
$ws.Range("F2").Value = "SMT_Line_8                    "
$ws.Range("G2").Value = 54
$ws.Range("I2").Value = "admin"
$ws.Range("J2").Value = "New entry (product exist)"

# Row 3
$ws.Range("A3").Value = 26
$ws.Range("B3").Value = "G5_00V6751_B_L8"
$ws.Range("C3").Value = "Bottom    "
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "SMT"
$ws.Range("F3").Value = "SMT_Line_9                    "
$ws.Range("G3").Value = 65
$ws.Range("I3").Value = "admin"
$ws.Range("J3").Value = "New entry  (Product does not exist)"

# Row 4
$ws.Range("A4").Value = 28
$ws.Range("B4").Value = "IB00VJ993_B"
$ws.Range("C4").Value = "Top       "
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = "SMT"
$ws.Range("F4").Value = "SMT_Line_8                    "
$ws.Range("G4").Value = 60
$ws.Range("I4").Value = "admin"
$ws.Range("J4").Value = "Duplicate entry from row 2 update Cycle time from 54>60"

# Row 5
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "IB00VJ993_B"
$ws.Range("C5").Value = "Top       "
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = "SMT"
$ws.Range("F5").Value = "SMT_Line_9                    "
$ws.Range("G5").Value = 40
$ws.Range("I5").Value = "admin"
$ws.Range("J5").Value = "Duplicate entry from row 2, different line"

# Row 6
$ws.Range("A6").Value = 27
$ws.Range("B6").Value = "G5_46X7155_B_L8"
$ws.Range("C6").Value = "Bottom    "
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = "SMT"
$ws.Range("F6").Value = "SMT_Line_10                    "
$ws.Range("G6").Value = 0
$ws.Range("I6").Value = "admin"
$ws.Range("J6").Value = "Invalid cycle time"

# Row 7
$ws.Range("A7").Value = 27
$ws.Range("B7").Value = "G5_46X7155_B_L8"
$ws.Range("C7").Value = "Bottom    "
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = "SMT"
$ws.Range("F7").Value = "SMT_Line_2345                    "
$ws.Range("G7").Value = 74
$ws.Range("I7").Value = "admin"
$ws.Range("J7").Value = "Invalid line"

# Row 8
$ws.Range("A8").Value = 23
$ws.Range("B8").Value = "gdrgdhdwetsfs"
$ws.Range("C8").Value = "fgjnfgfg"
$ws.Range("D8").Value = 6
$ws.Range("E8").Value = "SMT"
$ws.Range("F8").Value = "21352wbdfb"
$ws.Range("G8").Value = 543
$ws.Range("I8").Value = "admin"
$ws.Range("J8").Value = "Rubbish entry"

# ---- New rows 9-15 ----

# Row 9
$ws.Range("A9").Value = 32
$ws.Range("B9").Value = "G6_00V6897_T_L8"
$ws.Range("C9").Value = "Top       "
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = "SMT"
$ws.Range("F9").Value = "SMT_Line_12                   "
$ws.Range("G9").Value = 76
$ws.Range("H9").Value2 = 43822
$ws.Range("H9").NumberFormat = $ws.Range("H8").NumberFormat
$ws.Range("I9").Value = "admin"

# Row 10
$ws.Range("A10").Value = 33
$ws.Range("B10").Value = "G6_00V6900_B_L8"
$ws.Range("C10").Value = "Bottom    "
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = "SMT"
$ws.Range("F10").Value = "SMT_Line_12                    "
$ws.Range("G10").Value = 76
$ws.Range("H10").Value2 = 43822
$ws.Range("H10").NumberFormat = $ws.Range("H8").NumberFormat
$ws.Range("I10").Value = "admin"

# Row 11
$ws.Range("A11").Value = 34
$ws.Range("B11").Value = "G6_00V6907_T_L8"
$ws.Range("C11").Value = "Top       "
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = "SMT"
$ws.Range("F11").Value = "SMT_Line_11                    "
$ws.Range("G11").Value = 53
$ws.Range("H11").Value2 = 43822
$ws.Range("H11").NumberFormat = $ws.Range("H8").NumberFormat
$ws.Range("I11").Value = "admin"

# Row 12
$ws.Range("A12").Value = 35
$ws.Range("B12").Value = "G6_00V6910_T_L8"
$ws.Range("C12").Value = "Top       "
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = "SMT"
$ws.Range("F12").Value = "SMT_Line_10                    "
$ws.Range("G12").Value = 82
$ws.Range("H12").Value2 = 43822
$ws.Range("H12").NumberFormat = $ws.Range("H8").NumberFormat
$ws.Range("I12").Value = "admin"

# Row 13
$ws.Range("A13").Value = 36
$ws.Range("B13").Value = "G7_00MJ518_B_L8"
$ws.Range("C13").Value = "Bottom    "
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = "SMT"
$ws.Range("F13").Value = "SMT_Line_9                    "
$ws.Range("G13").Value = 63
$ws.Range("H13").Value2 = 43822
$ws.Range("H13").NumberFormat = $ws.Range("H8").NumberFormat
$ws.Range("I13").Value = "admin"

# Row 14
$ws.Range("A14").Value = 37
$ws.Range("B14").Value = "G7_00MJ521_B_L8"
$ws.Range("C14").Value = "Bottom    "
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = "SMT"
$ws.Range("F14").Value = "SMT_Line_10                    "
$ws.Range("G14").Value = 63
$ws.Range("H14").Value2 = 43822
$ws.Range("H14").NumberFormat = $ws.Range("H8").NumberFormat
$ws.Range("I14").Value = "admin"

# Row 15
$ws.Range("A15").Value = 38
$ws.Range("B15").Value = "G8_00VJ679_B_L8"
$ws.Range("C15").Value = "Bottom    "
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = "SMT"
$ws.Range("F15").Value = "SMT_Line_12                    "
$ws.Range("G15").Value = 63
$ws.Range("H15").Value2 = 43822
$ws.Range("H15").NumberFormat = $ws.Range("H8").NumberFormat
$ws.Range("I15").Value = "admin"

# ---- Column widths ----
$ws.Columns.Item(4).ColumnWidth = 6.1328125
$ws.Columns.Item(6).ColumnWidth = 14.53125
$ws.Columns.Item(7).ColumnWidth = 8.53125

# ---- Selection ----
$ws.Range("J11").Select()
